# Building effect size models
#
# Fill the "Study" label down column B for every study-group block on the
# "Temporal Normalization Studies" sheet: each group starts with a header
# row (e.g. row 2) that already carries the study name (with its section
# formatting) in column B, while the following "site/year" detail rows in
# that group had an empty, plainly-styled B cell (or, for the very last
# group, no B cell at all). Copy the header cell's value + formatting down
# across each of those detail rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Temporal Normalization Studies")

# Each entry: header row whose B cell holds the study name/style, and the
# row range of detail rows directly below it that belong to the same group.
$groups = @(
    @{ Header = 2;  First = 3;  Last = 4  },
    @{ Header = 5;  First = 6;  Last = 9  },
    @{ Header = 10; First = 11; Last = 20 },
    @{ Header = 21; First = 22; Last = 23 },
    @{ Header = 24; First = 25; Last = 27 },
    @{ Header = 28; First = 29; Last = 36 },
    @{ Header = 37; First = 38; Last = 41 },
    @{ Header = 42; First = 43; Last = 47 },
    @{ Header = 48; First = 49; Last = 53 },
    @{ Header = 54; First = 55; Last = 58 },
    @{ Header = 59; First = 60; Last = 61 },
    @{ Header = 62; First = 63; Last = 66 },
    @{ Header = 67; First = 68; Last = 69 },
    @{ Header = 70; First = 71; Last = 73 },
    @{ Header = 74; First = 75; Last = 77 },
    @{ Header = 78; First = 79; Last = 81 }
)

foreach ($g in $groups) {
    $headerCell = $ws.Cells.Item($g.Header, 2)
    $targetRange = $ws.Range($ws.Cells.Item($g.First, 2), $ws.Cells.Item($g.Last, 2))

    # Copy the header's formatting (style) down onto the detail rows, then
    # stamp the same study-name text into every cell of that range.
    $headerCell.Copy() | Out-Null
    $targetRange.PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    $targetRange.Value = $headerCell.Value2
}

$excel.CutCopyMode = $false

# Restore the sheet's last on-screen selection to match the edited area.
$ws.Activate()
$ws.Range("B78:B81").Select() | Out-Null

Write-Output "Filled study name down $($groups.Count) groups."
